# Auto-generated Excel COM-interop script
# Applies literal market-data value updates to the "Hades_Profits" workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1013.1667
$ws.Range("I18").Value = 1013.1667
$ws.Range("K18").Value = 1013.1667
$ws.Range("M18").Value = -729.1667
$ws.Range("H40").Value = 1750.25
$ws.Range("I40").Value = 3001
$ws.Range("J40").Value = 1333.3334
$ws.Range("K40").Value = 3001
$ws.Range("L40").Value = 1333.3334
$ws.Range("M40").Value = -2826
$ws.Range("N40").Value = -1683.3334
$ws.Range("H64").Value = 4787.5
$ws.Range("I64").Value = 3725
$ws.Range("J64").Value = 5850
$ws.Range("K64").Value = 3725
$ws.Range("L64").Value = 5850
$ws.Range("M64").Value = -3477
$ws.Range("N64").Value = -6346
$ws.Range("H67").Value = 4787.5
$ws.Range("I67").Value = 3725
$ws.Range("J67").Value = 5850
$ws.Range("K67").Value = 3725
$ws.Range("L67").Value = 5850
$ws.Range("M67").Value = -2867
$ws.Range("N67").Value = -7566
$ws.Range("H70").Value = 3236.8667
$ws.Range("I70").Value = 3571.4285
$ws.Range("J70").Value = 2944.125
$ws.Range("K70").Value = 10714.2855
$ws.Range("L70").Value = 8832.375
$ws.Range("M70").Value = -10444.2855
$ws.Range("N70").Value = -9372.375
$ws.Range("H73").Value = 3236.8667
$ws.Range("I73").Value = 3571.4285
$ws.Range("J73").Value = 2944.125
$ws.Range("K73").Value = 10714.2855
$ws.Range("L73").Value = 8832.375
$ws.Range("M73").Value = -9778.2855
$ws.Range("N73").Value = -10704.375
$ws.Range("H74").Value = 3919.439
$ws.Range("I74").Value = 3923.4412
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 3923.4412
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -2987.4412
$ws.Range("N74").Value = -5772
$ws.Range("H77").Value = 3919.439
$ws.Range("I77").Value = 3923.4412
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 19617.206
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -14937.206
$ws.Range("N77").Value = -28860
$ws.Range("H100").Value = 2298.4285
$ws.Range("I100").Value = 1076
$ws.Range("K100").Value = 1076
$ws.Range("M100").Value = -535
$ws.Range("H113").Value = 3432.4443
$ws.Range("I113").Value = 3278.3333
$ws.Range("K113").Value = 3278.3333
$ws.Range("M113").Value = -24.33329999999978
$ws.Range("H135").Value = 19581.555
$ws.Range("I135").Value = 23337.318
$ws.Range("J135").Value = 3056.2
$ws.Range("K135").Value = 210035.862
$ws.Range("L135").Value = 27505.8
$ws.Range("M135").Value = -207500.862
$ws.Range("N135").Value = -32575.8
$ws.Range("H137").Value = 2223790.2
$ws.Range("I137").Value = 2704091
$ws.Range("J137").Value = 2399.5
$ws.Range("K137").Value = 8112273
$ws.Range("L137").Value = 7198.5
$ws.Range("M137").Value = -8109723
$ws.Range("N137").Value = -12298.5
$ws.Range("H138").Value = 2720701.8
$ws.Range("I138").Value = 258166
$ws.Range("K138").Value = 774498
$ws.Range("M138").Value = -769358

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3473629
$ws.Range("I97").Value = 4465854.5
$ws.Range("J97").Value = 840.25
$ws.Range("K97").Value = 4465854.5
$ws.Range("L97").Value = 840.25
$ws.Range("M97").Value = -4465358.5
$ws.Range("N97").Value = -1832.25
$ws.Range("H102").Value = 20410464
$ws.Range("I102").Value = 35716060
$ws.Range("K102").Value = 35716060
$ws.Range("M102").Value = -35714438
$ws.Range("H132").Value = 101220.52
$ws.Range("I132").Value = 63852
$ws.Range("K132").Value = 191556
$ws.Range("M132").Value = -189026

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1116.1666
$ws.Range("I99").Value = 1149.1923
$ws.Range("J99").Value = 1062.5
$ws.Range("K99").Value = 1149.1923
$ws.Range("L99").Value = 1062.5
$ws.Range("M99").Value = 348.8077000000001
$ws.Range("N99").Value = -4058.5
$ws.Range("H105").Value = 33335658
$ws.Range("I105").Value = 62502584
$ws.Range("J105").Value = 2028.5714
$ws.Range("K105").Value = 62502584
$ws.Range("L105").Value = 2028.5714
$ws.Range("M105").Value = -62500837
$ws.Range("N105").Value = -5522.5714
$ws.Range("H134").Value = 5996.355
$ws.Range("I134").Value = 5495.654
$ws.Range("K134").Value = 16486.962
$ws.Range("M134").Value = -13951.962

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 28572666
$ws.Range("I58").Value = 38462524
$ws.Range("J58").Value = 1968.2222
$ws.Range("K58").Value = 38462524
$ws.Range("L58").Value = 1968.2222
$ws.Range("M58").Value = -38462321
$ws.Range("N58").Value = -2374.2222
$ws.Range("H62").Value = 2985.7144
$ws.Range("H65").Value = 2985.7144
$ws.Range("H105").Value = 1011.0526
$ws.Range("I105").Value = 1035.6428
$ws.Range("J105").Value = 942.2
$ws.Range("K105").Value = 1035.6428
$ws.Range("L105").Value = 942.2
$ws.Range("M105").Value = 711.3571999999999
$ws.Range("N105").Value = -4436.2
$ws.Range("H134").Value = 31347.36
$ws.Range("I134").Value = 3065.6785
$ws.Range("J134").Value = 103337.09
$ws.Range("K134").Value = 9197.0355
$ws.Range("L134").Value = 310011.27
$ws.Range("M134").Value = -6662.0355
$ws.Range("N134").Value = -315081.27
$ws.Range("H136").Value = 28572666
$ws.Range("I136").Value = 38462524
$ws.Range("J136").Value = 1968.2222
$ws.Range("K136").Value = 115387572
$ws.Range("L136").Value = 5904.6666
$ws.Range("M136").Value = -115385022
$ws.Range("N136").Value = -11004.6666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 936.3570999999999
$ws.Range("I34").Value = 443.7143
$ws.Range("J34").Value = 1429
$ws.Range("K34").Value = 1331.1429
$ws.Range("L34").Value = 4287
$ws.Range("M34").Value = -1247.1429
$ws.Range("N34").Value = -4455
$ws.Range("H39").Value = 2000
$ws.Range("J39").Value = 2750
$ws.Range("L39").Value = 8250
$ws.Range("N39").Value = -8838

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10560
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("N50").Value = -10996
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26372
$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81864
$ws.Range("H80").Value = 3952.5
$ws.Range("I80").Value = 3780
$ws.Range("J80").Value = 3961.5789
$ws.Range("K80").Value = 3780
$ws.Range("L80").Value = 3961.5789
$ws.Range("M80").Value = -2782
$ws.Range("N80").Value = -5957.5789
$ws.Range("H83").Value = 3952.5
$ws.Range("I83").Value = 3780
$ws.Range("J83").Value = 3961.5789
$ws.Range("K83").Value = 18900
$ws.Range("L83").Value = 19807.8945
$ws.Range("M83").Value = -13908
$ws.Range("N83").Value = -29791.8945
$ws.Range("H97").Value = 2072.8572
$ws.Range("I97").Value = 2218.3333
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 2218.3333
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -1722.3333
$ws.Range("N97").Value = -2192
$ws.Range("H132").Value = 68826.664
$ws.Range("I132").Value = 45373.39
$ws.Range("K132").Value = 136120.17
$ws.Range("M132").Value = -133590.17

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3651.8147
$ws.Range("I16").Value = 1719.4117
$ws.Range("J16").Value = 6936.9
$ws.Range("K16").Value = 1719.4117
$ws.Range("L16").Value = 6936.9
$ws.Range("M16").Value = -1549.4117
$ws.Range("N16").Value = -7276.9
$ws.Range("H68").Value = 1720
$ws.Range("I68").Value = 1533.3334
$ws.Range("K68").Value = 1533.3334
$ws.Range("M68").Value = -784.3334
$ws.Range("H71").Value = 1720
$ws.Range("I71").Value = 1533.3334
$ws.Range("K71").Value = 7666.666999999999
$ws.Range("M71").Value = -3922.666999999999
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 1360.0286
$ws.Range("I100").Value = 1128.619
$ws.Range("K100").Value = 1128.619
$ws.Range("M100").Value = -587.6189999999999
$ws.Range("H136").Value = 70695.80499999999
$ws.Range("I136").Value = 45390.8
$ws.Range("J136").Value = 176133.33
$ws.Range("K136").Value = 136172.4
$ws.Range("L136").Value = 528399.99
$ws.Range("M136").Value = -133622.4
$ws.Range("N136").Value = -533499.99

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 71433070
$ws.Range("I62").Value = 125003624
$ws.Range("J62").Value = 5668.6665
$ws.Range("K62").Value = 125003624
$ws.Range("L62").Value = 5668.6665
$ws.Range("M62").Value = -125003000
$ws.Range("N62").Value = -6916.6665
$ws.Range("H65").Value = 71433070
$ws.Range("I65").Value = 125003624
$ws.Range("J65").Value = 5668.6665
$ws.Range("K65").Value = 625018120
$ws.Range("L65").Value = 28343.3325
$ws.Range("M65").Value = -625015000
$ws.Range("N65").Value = -34583.3325
$ws.Range("H81").Value = 2835.9285
$ws.Range("I81").Value = 3600.5
$ws.Range("J81").Value = 2530.1
$ws.Range("K81").Value = 7201
$ws.Range("L81").Value = 5060.2
$ws.Range("M81").Value = -6140
$ws.Range("N81").Value = -7182.2
$ws.Range("H84").Value = 2835.9285
$ws.Range("I84").Value = 3600.5
$ws.Range("J84").Value = 2530.1
$ws.Range("K84").Value = 36005
$ws.Range("L84").Value = 25301
$ws.Range("M84").Value = -30701
$ws.Range("N84").Value = -35909
$ws.Range("H96").Value = 1586.1428
$ws.Range("I96").Value = 1586.1428
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1586.1428
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -213.1428000000001
$ws.Range("N96").ClearContents()
$ws.Range("H141").Value = 68732.5
$ws.Range("J141").Value = 68732.5
$ws.Range("L141").Value = 68732.5
$ws.Range("N141").Value = -79092.5
